# Add a new bullet "Post processing" right after the
# "Object Dissolve (The faster speed the faster dissolve)" bullet,
# inheriting that bullet's list/paragraph formatting (Listenabsatz,
# numId 5, Arial, en-US).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "Object Dissolve (The faster speed the faster dissolve)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the matched text, insert a sibling paragraph
    # after it (this copies the paragraph/run formatting of the source
    # paragraph), then move into that new paragraph and set its text.
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(4, 1) | Out-Null
    $rng.Text = "Post processing"
}
